# Empresas_Permisionarias.xlsx — "Add files via upload"
#
# The author re-uploaded the workbook for a new reporting period. The only
# substantive change is that the sole worksheet, previously labelled for
# cut "C_11", is relabelled "C_15.1" (the next table number in the series).
# Excel automatically keeps every reference to the old sheet name in sync
# when a sheet is renamed through the object model:
#   - the sheet tab itself
#   - the hidden _xlnm._FilterDatabase defined name that scopes the
#     existing AutoFilter on row 5 (B5) to this sheet
#   - the workbook's TitlesOfParts document property
# so a single `.Name` assignment on the worksheet reproduces all of it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "C_15.1"
